# "done surat pengajuan cuti" - turns a few literal placeholders in the
# surat-izin-penelitian template into merge fields:
#   1) Peneliti signature line (2nd table, row 1, col 2, last paragraph):
#      add a 3114-twip tab stop, re-indent the line (422 -> 2084 twips),
#      and put ${nama_mhs} on the underlined signature blank (previously
#      just a bare tab character).
#   2) "Koord. Program Studi ………….": the ellipsis run becomes ${prodi}.
#   3) Both "Karawang, 16 Februari 2024" dates: "16 Februari 2024"
#      becomes ${created_at}.
#
# Implementation notes (both worked around empirically against this
# runtime):
#   - ParagraphFormat.LeftIndent / TabStops mutations are no-ops for any
#     paragraph that lives inside a table here, so edit (1) is done by
#     splicing the paragraph's own OOXML back in via Range.InsertXML
#     (which *does* persist) instead of the indent/tab-stop properties.
#   - Touching the Tables/Cell collections (Tables.Item/.Cell(...)) then
#     corrupts subsequent Paragraphs.Item(n) lookups on this document (it
#     starts returning the wrong paragraph for every index), so the
#     target paragraphs below are located purely via $d.Paragraphs.Item,
#     using indices confirmed ahead of time against this document.
#   - Find.Execute run against the whole-document $d.Content range will
#     happily replace the (short) ellipsis/date substrings wherever they
#     occur as a *substring* of other, longer placeholder runs elsewhere
#     in the doc, so each Find below is scoped to just the one target
#     paragraph's Range.

$d = $word.ActiveDocument

# --- 1) Peneliti signature paragraph --------------------------------------
$sigPara = $d.Paragraphs.Item(75)

$newParaXml = '<w:p w14:paraId="67A13FF9" w14:textId="77777777" w:rsidR="00F53D08" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="TableParagraph"/><w:tabs><w:tab w:val="left" w:pos="3114"/><w:tab w:val="left" w:pos="4190"/></w:tabs><w:spacing w:line="210" w:lineRule="exact"/><w:ind w:left="2084"/><w:rPr><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:spacing w:val="-10"/><w:sz w:val="20"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:u w:val="single"/></w:rPr><w:t>${nama_mhs}</w:t></w:r><w:r><w:rPr><w:spacing w:val="-10"/><w:sz w:val="20"/></w:rPr><w:t>)</w:t></w:r></w:p>'
$sigPara.Range.InsertXML($newParaXml) | Out-Null

# --- 2) Koord. Program Studi ellipsis placeholder --------------------------
$prodiPara = $d.Paragraphs.Item(84)
$prodiPara.Range.Find.Execute("………….", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "`${prodi}", 2) | Out-Null

# --- 3) Both "16 Februari 2024" dates -> ${created_at} ---------------------
$datePara1 = $d.Paragraphs.Item(167)
$datePara1.Range.Find.Execute("16 Februari 2024", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "`${created_at}", 2) | Out-Null

$datePara2 = $d.Paragraphs.Item(206)
$datePara2.Range.Find.Execute("16 Februari 2024", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "`${created_at}", 2) | Out-Null
